# Update the "想去人数" (want-to-go count) column F values for the
# 展览 sheet and the 全部类型 sheet (which mirrors the same events,
# plus an extra row), to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Row -> new F value, as it applies to the "展览" sheet (sheet1).
$updatesSheet1 = @{
    2  = 1106
    3  = 421
    4  = 1513
    5  = 8785
    7  = 494
    9  = 294
    10 = 160
    11 = 22
    12 = 13
    13 = 3638
    15 = 369
    17 = 2418
    19 = 1125
    21 = 211
    22 = 2431
    23 = 70
}

# Row -> new F value, as it applies to the "全部类型" sheet (sheet4),
# which has an extra row (row 23) inserted relative to "展览", shifting
# the last update down to row 24.
$updatesSheet4 = @{
    2  = 1106
    3  = 421
    4  = 1513
    5  = 8785
    7  = 494
    9  = 294
    10 = 160
    11 = 22
    12 = 13
    13 = 3638
    15 = 369
    17 = 2418
    19 = 1125
    21 = 211
    22 = 2431
    24 = 70
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesSheet1.Keys) {
    $ws1.Range("F$row").Value = $updatesSheet1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesSheet4.Keys) {
    $ws4.Range("F$row").Value = $updatesSheet4[$row]
}
